# Regenerate save_data: write new "K" (Strike#) values into column G
# for the rows whose computed value changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new value for column G ("K")
$gUpdates = @{
    2  = 1
    3  = 3
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 3
    9  = 2
    10 = 2
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    16 = 1
    17 = 1
    18 = 3
    20 = 2
    21 = 1
    22 = 2
}

foreach ($row in $gUpdates.Keys) {
    $ws.Range("G$row").Value = $gUpdates[$row]
}
